# SARAALERT-1260: Allow vaccine table to be populated on import
# Adds two "Vaccine" blocks (Group Name, Product Name, Administration Date,
# Dose Number, Notes) to the Monitorees export template, columns CY:DH,
# with a few sample rows of demo data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (row 1) -- CY1:DH1
# ---------------------------------------------------------------------
$headers = @(
  "Vaccine 1 Group Name",
  "Vaccine 1 Product Name",
  "Vaccine 1 Administration Date",
  "Vaccine 1 Dose Number",
  "Vaccine 1 Notes",
  "Vaccine 2 Group Name",
  "Vaccine 2 Product Name",
  "Vaccine 2 Administration Date",
  "Vaccine 2 Dose Number",
  "Vaccine 2 Notes"
)

$startCol = 103  # CY
for ($i = 0; $i -lt $headers.Length; $i++) {
  $col = $startCol + $i
  $cell = $ws.Cells.Item(1, $col)
  if (($i -eq 2) -or ($i -eq 7)) {
    # Administration Date header columns (DA / DF) carry the Text style
    $cell.NumberFormat = "@"
  }
  $cell.Value = $headers[$i]
}

# ---------------------------------------------------------------------
# Data rows 2-7 -- CY:DH
# Columns within each block: Group(CY/DD) Product(CZ/DE) Date(DA/DF)
#                             Dose(DB/DG) Notes(DC/DH)
# ---------------------------------------------------------------------

# Row 2
$ws.Cells.Item(2, 103).Value = "COVID-19"
$ws.Cells.Item(2, 104).Value = "Moderna COVID-19 Vaccine"
$ws.Cells.Item(2, 105).NumberFormat = "@"
$ws.Cells.Item(2, 105).Value = "2020-06-01"
$ws.Cells.Item(2, 106).Value = 1
$ws.Cells.Item(2, 107).Value = "notes 1"
$ws.Cells.Item(2, 108).Value = "COVID-19"
$ws.Cells.Item(2, 109).Value = "Moderna COVID-19 Vaccine"
$ws.Cells.Item(2, 110).NumberFormat = "@"
$ws.Cells.Item(2, 110).Value = "2020-06-20"
$ws.Cells.Item(2, 111).Value = 2
$ws.Cells.Item(2, 112).Value = "notes 2"

# Row 3
$ws.Cells.Item(3, 103).Value = "COVID-19"
$ws.Cells.Item(3, 104).Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Cells.Item(3, 105).NumberFormat = "@"
$ws.Cells.Item(3, 105).Value = "2020-06-02"
$ws.Cells.Item(3, 106).Value = 1
$ws.Cells.Item(3, 108).Value = "COVID-19"
$ws.Cells.Item(3, 109).Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Cells.Item(3, 110).NumberFormat = "@"
$ws.Cells.Item(3, 110).Value = "2020-06-21"
$ws.Cells.Item(3, 111).Value = 2

# Row 4
$ws.Cells.Item(4, 103).Value = "COVID-19"
$ws.Cells.Item(4, 104).Value = "Unknown"
$ws.Cells.Item(4, 105).NumberFormat = "@"
$ws.Cells.Item(4, 105).Value = "2020-06-04"
$ws.Cells.Item(4, 106).Value = 1
$ws.Cells.Item(4, 108).Value = "COVID-19"
$ws.Cells.Item(4, 109).Value = "Unknown"
$ws.Cells.Item(4, 110).NumberFormat = "@"
$ws.Cells.Item(4, 110).Value = "2020-06-22"
$ws.Cells.Item(4, 111).Value = 2

# Row 5
$ws.Cells.Item(5, 103).Value = "COVID-19"
$ws.Cells.Item(5, 104).Value = "Moderna COVID-19 Vaccine"
$ws.Cells.Item(5, 105).NumberFormat = "@"
$ws.Cells.Item(5, 105).Value = "2020-06-01"
$ws.Cells.Item(5, 106).Value = 1

# Row 6
$ws.Cells.Item(6, 103).Value = "COVID-19"
$ws.Cells.Item(6, 104).Value = "Janssen (J&J) COVID-19 Vaccine"
$ws.Cells.Item(6, 105).NumberFormat = "@"
$ws.Cells.Item(6, 105).Value = "2020-06-03"
$ws.Cells.Item(6, 106).Value = 1

# Row 7
$ws.Cells.Item(7, 103).Value = "COVID-19"
$ws.Cells.Item(7, 104).Value = "Unknown"
$ws.Cells.Item(7, 105).NumberFormat = "@"
$ws.Cells.Item(7, 105).Value = "2020-06-02"
$ws.Cells.Item(7, 106).Value = 1

# ---------------------------------------------------------------------
# Column widths for the new columns (best effort -- engine quantizes to
# its own internal character grid, so these land as close as possible
# to the widths Excel originally wrote).
# ---------------------------------------------------------------------
$offset = 5 / 7
$ws.Columns.Item(103).ColumnWidth = (20.33203125 - $offset)
$ws.Columns.Item(104).ColumnWidth = (31 - $offset)
$ws.Columns.Item(105).ColumnWidth = (25.6640625 - $offset)
$ws.Columns.Item(106).ColumnWidth = (21.1640625 - $offset)
$ws.Columns.Item(107).ColumnWidth = (14.5 - $offset)
$ws.Columns.Item(108).ColumnWidth = (20.33203125 - $offset)
$ws.Columns.Item(109).ColumnWidth = (31 - $offset)
$ws.Columns.Item(110).ColumnWidth = (25.6640625 - $offset)
$ws.Columns.Item(111).ColumnWidth = (21.1640625 - $offset)
$ws.Columns.Item(112).ColumnWidth = (14.5 - $offset)

# ---------------------------------------------------------------------
# Reset selection back to A1 (the sheet had scrolled to CT1 / selected
# CZ10 before this edit).
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
